# Updated cryptos list with latest Price / Volume(1h) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '27.814.29'
Set-TextValue $ws.Range("E2") '  -0.49%  '
Set-TextValue $ws.Range("D3") '1.625.85'
Set-TextValue $ws.Range("E3") '  -0.54%  '
Set-TextValue $ws.Range("D4") '0.995'
Set-TextValue $ws.Range("E4") '  -0.27%  '
Set-TextValue $ws.Range("D5") '210.81'
Set-TextValue $ws.Range("E5") '  -0.55%  '
Set-TextValue $ws.Range("E6") '  -1.00%  '
Set-TextValue $ws.Range("E7") '  -0.23%  '
Set-TextValue $ws.Range("E8") '  -0.90%  '
Set-TextValue $ws.Range("D9") '0.256'
Set-TextValue $ws.Range("E9") '  -0.98%  '
Set-TextValue $ws.Range("E10") '  -1.14%  '
Set-TextValue $ws.Range("E11") '  -0.40%  '
Set-TextValue $ws.Range("D12") '1.856.22'
Set-TextValue $ws.Range("E12") '  -0.47%  '
Set-TextValue $ws.Range("D13") '1.625.52'
Set-TextValue $ws.Range("E13") '  -0.36%  '
Set-TextValue $ws.Range("E14") '  -1.33%  '
Set-TextValue $ws.Range("E15") '  -1.12%  '
Set-TextValue $ws.Range("D16") '64.85'
Set-TextValue $ws.Range("E16") '  -1.13%  '
Set-TextValue $ws.Range("D17") '27.824.57'
Set-TextValue $ws.Range("E17") '  -0.44%  '
Set-TextValue $ws.Range("D18") '227.87'
Set-TextValue $ws.Range("E18") '  -1.90%  '
Set-TextValue $ws.Range("E19") '  -1.16%  '
Set-TextValue $ws.Range("D20") '7.57'
Set-TextValue $ws.Range("E20") '  +0.27%  '
Set-TextValue $ws.Range("E21") '  -0.24%  '
Set-TextValue $ws.Range("D22") '4.33'
Set-TextValue $ws.Range("E22") '  -0.66%  '
Set-TextValue $ws.Range("D23") '9.92'
Set-TextValue $ws.Range("E23") '  -4.52%  '
Set-TextValue $ws.Range("E24") '  -0.55%  '
Set-TextValue $ws.Range("D25") '154.98'
Set-TextValue $ws.Range("E25") '  +0.31%  '
Set-TextValue $ws.Range("D27") '0.110'
Set-TextValue $ws.Range("E27") '  -0.75%  '
Set-TextValue $ws.Range("D28") '15.42'
Set-TextValue $ws.Range("E28") '  -1.54%  '
Set-TextValue $ws.Range("D29") '0.996'
Set-TextValue $ws.Range("E29") '  -0.30%  '
Set-TextValue $ws.Range("D30") '1.18'
Set-TextValue $ws.Range("E30") '  -0.25%  '
Set-TextValue $ws.Range("E31") '  -0.44%  '
Set-TextValue $ws.Range("E32") '  -0.25%  '
Set-TextValue $ws.Range("E33") '  +0.18%  '
Set-TextValue $ws.Range("D34") '1.403.71'
Set-TextValue $ws.Range("E34") '  -0.45%  '
Set-TextValue $ws.Range("D35") '1.60'
Set-TextValue $ws.Range("E35") '  +1.86%  '
Set-TextValue $ws.Range("E36") '  -0.69%  '
Set-TextValue $ws.Range("E37") '  -1.46%  '
Set-TextValue $ws.Range("E38") '  -0.90%  '
Set-TextValue $ws.Range("E39") '  -0.97%  '
Set-TextValue $ws.Range("E40") '  -2.97%  '
Set-TextValue $ws.Range("D41") '0.996'
Set-TextValue $ws.Range("E41") '  -0.21%  '
Set-TextValue $ws.Range("D42") '0.998'
Set-TextValue $ws.Range("E42") '  -2.43%  '
Set-TextValue $ws.Range("D43") '65.66'
Set-TextValue $ws.Range("E44") '  -0.85%  '
Set-TextValue $ws.Range("E45") '  -1.20%  '
Set-TextValue $ws.Range("D46") '1.766.04'
Set-TextValue $ws.Range("E46") '  -0.51%  '
Set-TextValue $ws.Range("E47") '  -4.23%  '
Set-TextValue $ws.Range("D48") '88.21'
Set-TextValue $ws.Range("E48") '  +0.04%  '
Set-TextValue $ws.Range("E49") '  +0.88%  '
Set-TextValue $ws.Range("E50") '  -0.65%  '
Set-TextValue $ws.Range("E51") '  +0.50%  '
